# Product Backlog v1 - apply commit "Add files via upload"
# Adds a Priority Key legend box, removes a stray note, reworks a few
# backlog item descriptions/ordering, restyles the "red" marker cell,
# and adds a totals formula at the bottom of the estimate column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the stray leftover comment in D10 ("this should be green?")
#    but keep its existing (yellow) fill style.
# ---------------------------------------------------------------------
$ws.Range("D10").ClearContents()

# ---------------------------------------------------------------------
# 2. Priority Key legend box in G4:J8
# ---------------------------------------------------------------------
$ws.Range("G4").Value = "Priority Key:"

$ws.Range("I5").Value = "Must Have"
$ws.Range("I5").Interior.Color = 5296274   # light green FF92D050

$ws.Range("I6").Value = "Could Have"
$ws.Range("I6").Interior.Color = 65535     # yellow FFFF00

$ws.Range("I7").Value = "Won't Have"
$ws.Range("I7").Interior.Color = 255       # red FF0000

# Column widths for the legend box columns (G, H, I)
$ws.Columns.Item(7).ColumnWidth = 11.42578125
$ws.Columns.Item(8).ColumnWidth = 11.140625
$ws.Columns.Item(9).ColumnWidth = 10.85546875

# Box border around H4:J8 - build it edge by edge (color first, then
# line style) so interior cells only pick up the shared edge they touch.
$topEdge = $ws.Range("H4:J4").Borders.Item(8)
$topEdge.Color = 0
$topEdge.LineStyle = 1

$leftEdge = $ws.Range("H4:H8").Borders.Item(7)
$leftEdge.Color = 0
$leftEdge.LineStyle = 1

$rightEdge = $ws.Range("J4:J8").Borders.Item(10)
$rightEdge.Color = 0
$rightEdge.LineStyle = 1

$bottomEdge = $ws.Range("H8:J8").Borders.Item(9)
$bottomEdge.Color = 0
$bottomEdge.LineStyle = 1

# ---------------------------------------------------------------------
# 3. Restyle the (empty) marker cell D26: red fill + yellow font
#    (previously it was a plain dark-red fill only).
# ---------------------------------------------------------------------
$ws.Range("D26").Interior.Color = 255      # red FF0000
$ws.Range("D26").Font.Color = 65535        # yellow FFFF00

# ---------------------------------------------------------------------
# 4. Rework the last few backlog items (reorder + reword).
# ---------------------------------------------------------------------
$ws.Range("B23").Value = "Implement recipe searching from ingredients"
$ws.Range("B24").Value = "Implement methods for the user to input the expiry dates of items"
$ws.Range("B25").Value = "Deduct items used in recipe from inventory"
$ws.Range("B26").Value = "Add game-ificiation elements like challenges and achievements"

# The custom row height follows the "inventory deduction" task, which
# moved from row 23 to row 25.
$ws.Rows.Item(23).RowHeight = 15
$ws.Rows.Item(23).EntireRow.AutoFit()
$ws.Rows.Item(25).RowHeight = 17.25

# ---------------------------------------------------------------------
# 5. Totals row: sum of estimated hours.
# ---------------------------------------------------------------------
$ws.Range("C28").Formula = "=SUM(C4:C26)"

# ---------------------------------------------------------------------
# 6. Selection / view bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("D10").Select()
